# Update "Luong ca nhan" report:
#  1. Fix "Nhóm dịch vụ" (G) values on sheet "Đơn 1 bác sĩ" rows 2 & 3 -> "Phun xăm"
#  2. Insert a new order row (row 4) on "Đơn 1 bác sĩ", pushing the "Tổng" row to 5
#     and refresh the totals on that row
#  3. Add a new sheet "Đơn thu nợ" with its header, one data row and a totals row

$wb = $excel.ActiveWorkbook

# ---------------------------------------------------------------------------
# 1 & 2: sheet "Đơn 1 bác sĩ"
# ---------------------------------------------------------------------------
$ws = $wb.Worksheets.Item("Đơn 1 bác sĩ")

$ws.Cells.Item(2, 7).Value = "Phun xăm"
$ws.Cells.Item(3, 7).Value = "Phun xăm"

# Push the existing "Tổng" row (row 4) down to row 5, freeing row 4 for the
# new order so the existing totals formatting / empty cells move with it.
$ws.Rows.Item(4).Insert()

$ws.Cells.Item(4, 1).Value = "HD-LUXURY"
$ws.Cells.Item(4, 2).Value = 561
# Force text so the "dd-mm-yyyy"-looking string is not auto-converted to a date serial
$ws.Cells.Item(4, 3).NumberFormat = "@"
$ws.Cells.Item(4, 3).Value = "07-16-2024"
$ws.Cells.Item(4, 4).Value = "SÓC TRĂNG"
$ws.Cells.Item(4, 5).Value = "diệp ngọc anh "
$ws.Cells.Item(4, 6).Value = "CTV"
$ws.Cells.Item(4, 7).Value = "Phun xăm"
$ws.Cells.Item(4, 8).Value = "Điêu khắc mày"
$ws.Cells.Item(4, 9).Value = 0
$ws.Cells.Item(4, 10).Value = 500000
$ws.Cells.Item(4, 11).Value = "Kha Như Huỳnh "
$ws.Cells.Item(4, 12).Value = 800000
$ws.Cells.Item(4, 13).Value = 1300000
$ws.Cells.Item(4, 14).Value = 1300000
$ws.Cells.Item(4, 15).Value = 0
$ws.Cells.Item(4, 16).Value = 1300000
$ws.Cells.Item(4, 17).Value = 0
$ws.Cells.Item(4, 18).Value = "Bác Sĩ Ngoài"
$ws.Cells.Item(4, 19).Value = 0
$ws.Cells.Item(4, 20).Value = 0
$ws.Cells.Item(4, 21).Value = 0
$ws.Cells.Item(4, 22).Value = 0
$ws.Cells.Item(4, 23).Value = 0
$ws.Cells.Item(4, 24).Value = 0
$ws.Cells.Item(4, 25).Value = 0.02
$ws.Cells.Item(4, 26).Value = 0
$ws.Cells.Item(4, 27).Value = 26000

# Refresh the totals that now live on row 5
$ws.Cells.Item(5, 2).Value = 3
$ws.Cells.Item(5, 10).Value = 7000000
$ws.Cells.Item(5, 12).Value = 800000
$ws.Cells.Item(5, 13).Value = 7800000
$ws.Cells.Item(5, 14).Value = 7800000
$ws.Cells.Item(5, 15).Value = 0
$ws.Cells.Item(5, 16).Value = 7800000
$ws.Cells.Item(5, 17).Value = 0
$ws.Cells.Item(5, 24).Value = 0.1
$ws.Cells.Item(5, 25).Value = 0.02
$ws.Cells.Item(5, 26).Value = 100000
$ws.Cells.Item(5, 27).Value = 26000

# ---------------------------------------------------------------------------
# 3: new sheet "Đơn thu nợ"
# ---------------------------------------------------------------------------
$lastSheet = $wb.Worksheets.Item($wb.Worksheets.Count)
$ws3 = $wb.Worksheets.Add($null, $lastSheet)
$ws3.Name = "Đơn thu nợ"

$headers = @(
    "Ngày thực hiện",
    "Ngày thu",
    "notion id",
    "Tiền tố",
    "Mã đơn thu nợ",
    "Cơ sở",
    "id đơn nợ",
    "Lượng thu",
    "Đơn nợ",
    "Nguồn khách",
    "Sale chính",
    "Đơn giá gốc",
    "Sale phụ",
    "Upsale",
    "Bác sĩ 1",
    "Bác sĩ 2",
    "Thanh toán lần đầu",
    "Đã thanh toán",
    "Tỉ lệ chiết khấu sale chính",
    "Tỉ lệ chiết khấu sale phụ",
    "id sale chính",
    "id sale phụ",
    "id bác sĩ 1",
    "id bác sĩ 2",
    "Chiết khấu bác sĩ 1",
    "Chiết khấu bác sĩ 2",
    "Chiết khấu sale chính",
    "Chiết khấu sale phụ"
)
for ($i = 0; $i -lt $headers.Count; $i++) {
    $ws3.Cells.Item(1, $i + 1).Value = $headers[$i]
}

$ws3.Cells.Item(2, 1).NumberFormat = "@"
$ws3.Cells.Item(2, 1).Value = "06-30-2024"
$ws3.Cells.Item(2, 2).NumberFormat = "@"
$ws3.Cells.Item(2, 2).Value = "07-07-2024"
$ws3.Cells.Item(2, 3).Value = "c050b03e-22d0-4cdc-8e50-69d05748efc3"
$ws3.Cells.Item(2, 4).Value = "TN"
$ws3.Cells.Item(2, 5).Value = 142
$ws3.Cells.Item(2, 6).Value = "SÓC TRĂNG"
$ws3.Cells.Item(2, 7).Value = "45183006-5603-42d6-a5bf-0b84f7b03e93"
$ws3.Cells.Item(2, 8).Value = 2000000
$ws3.Cells.Item(2, 9).Value = "HD-LUXURY-498"
$ws3.Cells.Item(2, 10).Value = "Cá nhân"
$ws3.Cells.Item(2, 11).Value = "Lê Đình Hậu"
$ws3.Cells.Item(2, 12).Value = 7000000
$ws3.Cells.Item(2, 13).Value = 0
$ws3.Cells.Item(2, 14).Value = 0
$ws3.Cells.Item(2, 15).Value = "Bác Sĩ Ngoài"
$ws3.Cells.Item(2, 16).Value = 0
$ws3.Cells.Item(2, 17).Value = 5000000
$ws3.Cells.Item(2, 18).Value = 7000000
$ws3.Cells.Item(2, 19).Value = 0.13
$ws3.Cells.Item(2, 20).Value = 0
$ws3.Cells.Item(2, 21).Value = "c463b1a9-4fb2-4258-87a7-44193ba02405"
$ws3.Cells.Item(2, 22).Value = 0
$ws3.Cells.Item(2, 23).Value = "545126c4-c319-4d90-a506-2627e3e232a0"
$ws3.Cells.Item(2, 24).Value = 0
$ws3.Cells.Item(2, 25).Value = 200000
$ws3.Cells.Item(2, 26).Value = 0
$ws3.Cells.Item(2, 27).Value = 260000
$ws3.Cells.Item(2, 28).Value = 0

$blankTextCols = @(1, 2, 3, 6, 7, 9, 10, 11, 13, 15, 16, 21, 22, 23, 24)
foreach ($col in $blankTextCols) {
    $ws3.Cells.Item(3, $col).NumberFormat = "General"
}

$ws3.Cells.Item(3, 4).Value = "Tổng"
$ws3.Cells.Item(3, 5).Value = 1
$ws3.Cells.Item(3, 8).Value = 2000000
$ws3.Cells.Item(3, 12).Value = 7000000
$ws3.Cells.Item(3, 14).Value = 0
$ws3.Cells.Item(3, 17).Value = 5000000
$ws3.Cells.Item(3, 18).Value = 7000000
$ws3.Cells.Item(3, 19).Value = 0.13
$ws3.Cells.Item(3, 20).Value = 0
$ws3.Cells.Item(3, 25).Value = 200000
$ws3.Cells.Item(3, 26).Value = 0
$ws3.Cells.Item(3, 27).Value = 260000
$ws3.Cells.Item(3, 28).Value = 0
